$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column, matching style of existing header cells (G1 etc.)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save flag values for each data row (2-35)
$saveValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 0
    35 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
